# Actualizacion automatica del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 347 had a pending result; fill in the now-known outcome ---
$ws.Cells.Item(347, 7).Value = "Fallo"
$ws.Cells.Item(347, 8).Value = -1

# --- Append the newly-tracked matches as rows 348-353 ---
$newRows = @(
    @{ Row = 348; A = 14428736; B = "2025-08-20"; C = "Yunchaokete Bu";   D = "Mariano Navone";           E = "Gana Mariano Navone";        F = 2.3 },
    @{ Row = 349; A = 14428728; B = "2025-08-20"; C = "Sebastián Báez";   D = "Botic Van de Zandschulp";  E = "Gana Sebastián Báez";        F = 2.5 },
    @{ Row = 350; A = 14427816; B = "2025-08-20"; C = "Rebecca Sramkova"; D = "Leylah Fernandez";         E = "Gana Rebecca Sramkova";      F = 3.5 },
    @{ Row = 351; A = 14427814; B = "2025-08-21"; C = "Linda Noskova";    D = "Tatjana Maria";            E = "Gana Tatjana Maria";         F = 4 },
    @{ Row = 352; A = 14427817; B = "2025-08-21"; C = "Emma Navarro";     D = "Alycia Parks";             E = "Gana Emma Navarro";          F = 1.36 },
    @{ Row = 353; A = 14427822; B = "2025-08-21"; C = "Ajla Tomljanovic"; D = "Ekaterina Alexandrova";    E = "Gana Ekaterina Alexandrova"; F = 1.33 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A

    # Write the date as literal text, not as an auto-converted date serial.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    # G (resultado) / H (profit) stay blank - result is still pending for these matches.
}
